$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Final_Matches"

$newSheet.Cells.Item(1, 1).Value = 'AZ.CT/LABEL'
$newSheet.Cells.Item(1, 2).Value = 'ASCTB.CT/LABEL'
$newSheet.Cells.Item(2, 1).Value = 'T cell'
$newSheet.Cells.Item(2, 2).Value = 'T cell'
$newSheet.Cells.Item(3, 1).Value = 'mast cell'
$newSheet.Cells.Item(3, 2).Value = 'mast cell'
$newSheet.Cells.Item(4, 1).Value = 'endothelial cell'
$newSheet.Cells.Item(4, 2).Value = 'endothlial cell'
$newSheet.Cells.Item(5, 1).Value = 'monocyte'
$newSheet.Cells.Item(5, 2).Value = 'monocyte'
$newSheet.Cells.Item(6, 1).Value = 'kidney granular cell'
$newSheet.Cells.Item(6, 2).Value = 'kidney granular cell'
$newSheet.Cells.Item(7, 1).Value = 'glomerular visceral epithelial cell'
$newSheet.Cells.Item(7, 2).Value = 'glomerular visceral epithelial cell'
$newSheet.Cells.Item(8, 1).Value = 'neutrophil'
$newSheet.Cells.Item(8, 2).Value = 'neutrophil'
$newSheet.Cells.Item(9, 1).Value = 'plasma cell'
$newSheet.Cells.Item(9, 2).Value = 'plasma cell'
$newSheet.Cells.Item(10, 1).Value = 'endothelial cell of lymphatic vessel'
$newSheet.Cells.Item(10, 2).Value = 'endothelial cell of lymphatic vessel'
$newSheet.Cells.Item(11, 1).Value = 'renal beta-intercalated cell'
$newSheet.Cells.Item(11, 2).Value = 'renal beta-intercalated cell'
$newSheet.Cells.Item(12, 1).Value = 'epithelial cell of proximal tubule'
$newSheet.Cells.Item(12, 2).Value = 'epithelial cell of proximal tubule'
$newSheet.Cells.Item(13, 1).Value = 'parietal epithelial cell'
$newSheet.Cells.Item(13, 2).Value = 'parietal epithelial cell'
$newSheet.Cells.Item(14, 1).Value = 'kidney outer medulla collecting duct principal cell'
$newSheet.Cells.Item(14, 2).Value = 'kidney outer medulla collecting duct principal cell'
$newSheet.Cells.Item(15, 1).Value = 'kidney outer medulla collecting duct intercalated cell'
$newSheet.Cells.Item(15, 2).Value = 'kidney outer medulla collecting duct intercalated cell'
$newSheet.Cells.Item(16, 1).Value = 'kidney connecting tubule epithelial cell'
$newSheet.Cells.Item(16, 2).Value = 'kidney connecting tubule epithelial cell'
$newSheet.Cells.Item(17, 1).Value = 'kidney distal convoluted tubule epithelial cell'
$newSheet.Cells.Item(17, 2).Value = 'kidney distal convoluted tubule epithelial cell'
$newSheet.Cells.Item(18, 1).Value = 'macula densa epithelial cell'
$newSheet.Cells.Item(18, 2).Value = 'macula densa epithelial cell'
$newSheet.Cells.Item(19, 1).Value = 'glomerular capillary endothelial cell'
$newSheet.Cells.Item(19, 2).Value = 'glomerular capillary endothelial cell'
$newSheet.Cells.Item(20, 1).Value = 'peritubular capillary endothelial cell'
$newSheet.Cells.Item(20, 2).Value = 'peritubular capillary endothelial cell'
$newSheet.Cells.Item(21, 1).Value = 'kidney afferent arteriole endothelial cell'
$newSheet.Cells.Item(21, 2).Value = 'kidney afferent arteriole endothelial cell'
$newSheet.Cells.Item(22, 1).Value = 'kidney loop of Henle thick ascending limb epithelial cell'
$newSheet.Cells.Item(22, 2).Value = 'kidney loop of Henle thick ascending limb epithelial cell'
$newSheet.Cells.Item(23, 1).Value = 'kidney loop of Henle thin ascending limb epithelial cell_x000D_' + "`n" + ''
$newSheet.Cells.Item(23, 2).Value = 'kidney loop of Henle thin ascending limb epithelial cell'
$newSheet.Cells.Item(24, 1).Value = 'kidney loop of Henle medullary thick ascending limb epithelial cell'
$newSheet.Cells.Item(24, 2).Value = 'kidney loop of Henle medullary thick ascending limb epithelial cell'
$newSheet.Cells.Item(25, 1).Value = 'kidney loop of Henle cortical thick ascending limb epithelial cell'
$newSheet.Cells.Item(25, 2).Value = 'kidney loop of Henle cortical thick ascending limb epithelial cell'
$newSheet.Cells.Item(26, 1).Value = 'kidney loop of Henle thin descending limb epithelial cell'
$newSheet.Cells.Item(26, 2).Value = 'kidney loop of Henle thin descending limb epithelial cell'
$newSheet.Cells.Item(27, 1).Value = 'vasa recta ascending limb cell_x000D_' + "`n" + ''
$newSheet.Cells.Item(27, 2).Value = 'vasa recta ascending limb cell'
$newSheet.Cells.Item(28, 1).Value = 'vasa recta descending limb cell'
$newSheet.Cells.Item(28, 2).Value = 'vasa recta descending limb cell'
$newSheet.Cells.Item(29, 1).Value = 'kidney collecting duct principal cell'
$newSheet.Cells.Item(29, 2).Value = 'kidney collecting duct principal cell'
$newSheet.Cells.Item(30, 1).Value = 'mature B cell'
$newSheet.Cells.Item(30, 2).Value = 'B cell'
$newSheet.Cells.Item(31, 1).Value = 'renal alpha-intercalated cell'
$newSheet.Cells.Item(31, 2).Value = 'kidney collecting duct intercalated cell'
$newSheet.Cells.Item(32, 1).Value = 'kidney collecting duct principal cell'
$newSheet.Cells.Item(32, 2).Value = 'kidney cortex collecting duct principal cell'
$newSheet.Cells.Item(33, 1).Value = 'kidney collecting duct principal cell'
$newSheet.Cells.Item(33, 2).Value = 'kidney inner medulla collecting duct principal cell'
$newSheet.Cells.Item(34, 1).Value = 'fibroblast'
$newSheet.Cells.Item(34, 2).Value = 'kidney interstitial fibroblast'
$newSheet.Cells.Item(35, 1).Value = 'macrophage'
$newSheet.Cells.Item(35, 2).Value = 'kidney resident macrophage'
$newSheet.Cells.Item(36, 1).Value = 'lymphocyte'
$newSheet.Cells.Item(36, 2).Value = 'natural killer cell'
$newSheet.Cells.Item(37, 1).Value = 'T cell'
$newSheet.Cells.Item(37, 2).Value = 'mature NK T cell'
$newSheet.Cells.Item(38, 1).Value = 'plasmacytoid dendritic cell'
$newSheet.Cells.Item(38, 2).Value = 'plasmacytoid dendritic cell, human'
$newSheet.Cells.Item(39, 1).Value = 'T cell'
$newSheet.Cells.Item(39, 2).Value = 'cytotoxic T cell'
$newSheet.Cells.Item(40, 1).Value = 'lymphocyte'
$newSheet.Cells.Item(40, 2).Value = 'B cell'

# Header formatting (bold, centered, bordered) to match existing header style
$header = $newSheet.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Column widths
$newSheet.Columns.Item(1).ColumnWidth = 59.88671875
$newSheet.Columns.Item(2).ColumnWidth = 59.88671875

# Selection matching target view state, and make this sheet active/selected (last edited)
$newSheet.Range("B28").Select()